# Update iServ stats for 2025-11 (row 24) to match the latest refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B24").Value = 6421
$ws.Range("C24").Value = 1004
$ws.Range("D24").Value = 5989557
$ws.Range("E24").Value = 932.8075066189067
$ws.Range("F24").Value = 9.461302420729623
$ws.Range("G24").Value = 4.041450777202082
$ws.Range("H24").Value = 26.88201768768701
